$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 112.333336
$ws.Range("I4").Value = 108.875
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 108.875
$ws.Range("L4").Value = 140
$ws.Range("M4").Value = 5.125
$ws.Range("N4").Value = -368
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H29").Value = 204
$ws.Range("I29").Value = 56
$ws.Range("K29").Value = 168
$ws.Range("M29").Value = 113
$ws.Range("H38").Value = 154
$ws.Range("I38").Value = 154
$ws.Range("K38").Value = 462
$ws.Range("M38").Value = -90
$ws.Range("H43").Value = 1657.1428
$ws.Range("J43").Value = 2150
$ws.Range("L43").Value = 2150
$ws.Range("N43").Value = -2288
$ws.Range("H82").Value = 246
$ws.Range("I82").Value = 246
$ws.Range("K82").Value = 738
$ws.Range("M82").Value = -332
$ws.Range("H85").Value = 246
$ws.Range("I85").Value = 246
$ws.Range("K85").Value = 738
$ws.Range("M85").Value = 666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2049.5
$ws.Range("I61").Value = 1914
$ws.Range("K61").Value = 1914
$ws.Range("M61").Value = -1702
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2049.923
$ws.Range("I132").Value = 2137.4167
$ws.Range("K132").Value = 6412.250100000001
$ws.Range("M132").Value = -3882.250100000001
$ws.Range("H136").Value = 2049.5
$ws.Range("I136").Value = 1914
$ws.Range("K136").Value = 5742
$ws.Range("M136").Value = -3192

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3675.6
$ws.Range("I20").Value = 3280
$ws.Range("K20").Value = 3280
$ws.Range("M20").Value = -3033
$ws.Range("H86").Value = 4155.5557
$ws.Range("I86").Value = 5475
$ws.Range("K86").Value = 5475
$ws.Range("M86").Value = -4352
$ws.Range("H89").Value = 4155.5557
$ws.Range("I89").Value = 5475
$ws.Range("K89").Value = 27375
$ws.Range("M89").Value = -21759
$ws.Range("H105").Value = 3356
$ws.Range("I105").Value = 3356
$ws.Range("K105").Value = 3356
$ws.Range("M105").Value = -1609
$ws.Range("H134").Value = 6890.6665
$ws.Range("I134").Value = 5668.684
$ws.Range("K134").Value = 17006.052
$ws.Range("M134").Value = -14471.052

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1440.6666
$ws.Range("J31").Value = 1664.6666
$ws.Range("L31").Value = 1664.6666
$ws.Range("N31").Value = -2254.6666
$ws.Range("H34").Value = 1440.6666
$ws.Range("J34").Value = 1664.6666
$ws.Range("L34").Value = 1664.6666
$ws.Range("N34").Value = -2068.6666
$ws.Range("H58").Value = 1615.1765
$ws.Range("I58").Value = 1615
$ws.Range("J58").Value = 1615.5
$ws.Range("K58").Value = 1615
$ws.Range("L58").Value = 1615.5
$ws.Range("M58").Value = -1412
$ws.Range("N58").Value = -2021.5
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H134").Value = 3830.9333
$ws.Range("I134").Value = 2708
$ws.Range("K134").Value = 8124
$ws.Range("M134").Value = -5589
$ws.Range("H136").Value = 1615.1765
$ws.Range("I136").Value = 1615
$ws.Range("J136").Value = 1615.5
$ws.Range("K136").Value = 4845
$ws.Range("L136").Value = 4846.5
$ws.Range("M136").Value = -2295
$ws.Range("N136").Value = -9946.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 733.6667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 733.6667
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2201.0001
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -2671.0001
$ws.Range("H80").Value = 7481
$ws.Range("J80").Value = 7804.8335
$ws.Range("L80").Value = 23414.5005
$ws.Range("N80").Value = -25286.5005
$ws.Range("H83").Value = 7481
$ws.Range("J83").Value = 7804.8335
$ws.Range("L83").Value = 70243.5015
$ws.Range("N83").Value = -79603.5015
$ws.Range("H134").Value = 143718.72
$ws.Range("I134").Value = 143718.72
$ws.Range("K134").Value = 431156.16
$ws.Range("M134").Value = -426086.16
$ws.Range("H140").Value = 2153.7778
$ws.Range("I140").Value = 1193.4286
$ws.Range("K140").Value = 3580.2858
$ws.Range("M140").Value = 1599.7142

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4374.75
$ws.Range("I36").Value = 4666.3335
$ws.Range("J36").Value = 3500
$ws.Range("K36").Value = 4666.3335
$ws.Range("L36").Value = 3500
$ws.Range("M36").Value = -4181.3335
$ws.Range("N36").Value = -4470
$ws.Range("H40").Value = 22000
$ws.Range("J40").Value = 22000
$ws.Range("L40").Value = 22000
$ws.Range("N40").Value = -22302
$ws.Range("H46").Value = 9800
$ws.Range("I46").Value = 9500
$ws.Range("K46").Value = 9500
$ws.Range("M46").Value = -9344

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 38000
$ws.Range("J50").Value = 38000
$ws.Range("L50").Value = 38000
$ws.Range("N50").Value = -39274
$ws.Range("H56").Value = 32000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H63").Value = 36361.668
$ws.Range("J63").Value = 36361.668
$ws.Range("L63").Value = 36361.668
$ws.Range("N63").Value = -37859.668
$ws.Range("H66").Value = 36361.668
$ws.Range("J66").Value = 36361.668
$ws.Range("L66").Value = 109085.004
$ws.Range("N66").Value = -116573.004
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("K74").Value = 25000
$ws.Range("M74").Value = -24002
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70008

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H51").Value = 54961.668
$ws.Range("I51").Value = 54770
$ws.Range("J51").Value = 55000
$ws.Range("K51").Value = 54770
$ws.Range("L51").Value = 55000
$ws.Range("M51").Value = -54260
$ws.Range("N51").Value = -56020
$ws.Range("H52").Value = 8550.75
$ws.Range("I52").Value = 8067.6665
$ws.Range("K52").Value = 8067.6665
$ws.Range("M52").Value = -7841.6665
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 4783.35
$ws.Range("I136").Value = 4937.933
$ws.Range("J136").Value = 4319.6
$ws.Range("K136").Value = 14813.799
$ws.Range("L136").Value = 12958.8
$ws.Range("M136").Value = -12263.799
$ws.Range("N136").Value = -18058.8
